$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D7").Value = "2016-03-09 14:02:48"
$wsZh.Range("D10").Value = "2016-03-09 14:02:48"
$wsZh.Range("D11").Value = "2016-03-09 14:02:48"
$wsZh.Range("D12").Value = "2016-03-09 14:02:48"
$wsZh.Range("D13").Value = "2016-03-09 14:02:48"
$wsZh.Range("D14").Value = "2016-03-09 14:02:48"
$wsZh.Range("D15").Value = "2016-03-09 14:02:48"
$wsZh.Range("D16").Value = "2016-03-09 14:02:48"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D7").Value = "2016-03-09 14:02:54"
$wsDe.Range("D10").Value = "2016-03-09 14:02:54"
$wsDe.Range("D11").Value = "2016-03-09 14:02:54"
$wsDe.Range("D12").Value = "2016-03-09 14:02:54"
$wsDe.Range("D13").Value = "2016-03-09 14:02:54"
$wsDe.Range("D14").Value = "2016-03-09 14:02:54"
$wsDe.Range("D15").Value = "2016-03-09 14:02:54"
$wsDe.Range("D16").Value = "2016-03-09 14:02:54"
